$p = $ppt.ActivePresentation

# Slide 8 ("EntityManagerFactory") -> shape "object 4" holds the Java code sample.
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Replace the email domain "pluralsight.com" -> "mycompany.com" (keeps the
# hyperlink run formatting; PowerPoint splits the run at the edited substring).
$fullText = $tr.Text
$needle = "pluralsight"
$idx = $fullText.IndexOf($needle)
if ($idx -ge 0) {
    $target = $tr.Characters($idx + 1, $needle.Length)
    $target.Text = "mycompany"
}

# The text box is an auto-fit shape; the extra characters push it onto one
# more wrapped line, so PowerPoint grows the shape's height to fit.
$sh.Height = 605.15
